# Fix some bugs about the story class.
# Rows 3 and 5 had their B:I data (Product Manager, Item, Jira Ticket, M109/M110,
# Regulatory/Compliance Urgency, Partner Priority, Product Value) swapped onto the
# wrong row. Swap them back so "John Rease / Testing / proj-003 / M109" data is on
# row 3 and "Tony Wei / Sign Off / proj-005 / M110" data is on row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "E", "G", "H", "I")

foreach ($col in $columns) {
    $cell3 = $ws.Range("$col`3")
    $cell5 = $ws.Range("$col`5")

    $val3 = $cell3.Value2
    $val5 = $cell5.Value2

    $cell3.Value2 = $val5
    $cell5.Value2 = $val3
}
